$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save") - copy the existing header style (from G1) so it
# reuses the same cellXf as the rest of the header row, then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cells H2 and H3 with numeric value 0 (plain, unstyled like B2:G3)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
